# Auto-generated script applying scraped price/profit data updates
# to the Cerberus_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 13891847
$ws.Range("I28").Value = 17546238
$ws.Range("K28").Value = 17546238
$ws.Range("M28").Value = -17545753
$ws.Range("H32").Value = 3499
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3499
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3499
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4151
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H69").Value = 11806.25
$ws.Range("I69").Value = 10832.571
$ws.Range("J69").Value = 13169.4
$ws.Range("K69").Value = 32497.713
$ws.Range("L69").Value = 39508.2
$ws.Range("M69").Value = -31623.713
$ws.Range("N69").Value = -41256.2
$ws.Range("H72").Value = 11806.25
$ws.Range("I72").Value = 10832.571
$ws.Range("J72").Value = 13169.4
$ws.Range("K72").Value = 97493.139
$ws.Range("L72").Value = 118524.6
$ws.Range("M72").Value = -93125.139
$ws.Range("N72").Value = -127260.6
$ws.Range("H86").Value = 4349.1113
$ws.Range("J86").Value = 4071
$ws.Range("L86").Value = 4071
$ws.Range("N86").Value = -6317
$ws.Range("H89").Value = 4349.1113
$ws.Range("J89").Value = 4071
$ws.Range("L89").Value = 20355
$ws.Range("N89").Value = -31587
$ws.Range("H92").Value = 499.66666
$ws.Range("I92").Value = 499.5
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 499.5
$ws.Range("L92").Value = 500
$ws.Range("M92").Value = 748.5
$ws.Range("N92").Value = -2996
$ws.Range("H93").Value = 53198.8
$ws.Range("J93").Value = 53198.8
$ws.Range("L93").Value = 53198.8
$ws.Range("N93").Value = -58190.8
$ws.Range("H98").Value = 7655.9473
$ws.Range("I98").Value = 1731.9231
$ws.Range("K98").Value = 1731.9231
$ws.Range("M98").Value = -233.9231
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H118").Value = 1741.5714
$ws.Range("I118").Value = 840.6
$ws.Range("K118").Value = 2521.8
$ws.Range("M118").Value = -864.8000000000002
$ws.Range("H122").Value = 7655.9473
$ws.Range("I122").Value = 1731.9231
$ws.Range("K122").Value = 5195.7693
$ws.Range("M122").Value = -2745.7693
$ws.Range("H125").Value = 5656.1665
$ws.Range("I125").Value = 4470
$ws.Range("K125").Value = 40230
$ws.Range("M125").Value = -37770
$ws.Range("H132").Value = 4146.2383
$ws.Range("I132").Value = 4101.049
$ws.Range("K132").Value = 12303.147
$ws.Range("M132").Value = -9773.147000000001
$ws.Range("H133").Value = 54912.145
$ws.Range("J133").Value = 54912.145
$ws.Range("L133").Value = 54912.145
$ws.Range("N133").Value = -65032.145
$ws.Range("H135").Value = 2915.2
$ws.Range("I135").Value = 2662.2083
$ws.Range("K135").Value = 23959.8747
$ws.Range("M135").Value = -21424.8747
$ws.Range("H137").Value = 1238.4445
$ws.Range("I137").Value = 859.6923
$ws.Range("K137").Value = 2579.0769
$ws.Range("M137").Value = -29.07690000000002
$ws.Range("H138").Value = 3220.7378
$ws.Range("I138").Value = 4369.6313
$ws.Range("J138").Value = 2701
$ws.Range("K138").Value = 13108.8939
$ws.Range("L138").Value = 8103
$ws.Range("M138").Value = -7968.893899999999
$ws.Range("N138").Value = -18383
$ws.Range("H141").Value = 11499
$ws.Range("I141").Value = 5999
$ws.Range("K141").Value = 17997
$ws.Range("M141").Value = -12817

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 848.7917
$ws.Range("I2").Value = 643.65
$ws.Range("K2").Value = 643.65
$ws.Range("M2").Value = -530.65
$ws.Range("H5").Value = 37.5
$ws.Range("I5").Value = 37.5
$ws.Range("K5").Value = 37.5
$ws.Range("M5").Value = 74.5
$ws.Range("H32").Value = 3941.275
$ws.Range("J32").Value = 14164
$ws.Range("L32").Value = 14164
$ws.Range("N32").Value = -14738
$ws.Range("H61").Value = 3649.6155
$ws.Range("I61").Value = 3299.6086
$ws.Range("J61").Value = 6333
$ws.Range("K61").Value = 3299.6086
$ws.Range("L61").Value = 6333
$ws.Range("M61").Value = -3087.6086
$ws.Range("N61").Value = -6757
$ws.Range("H74").Value = 1687.6471
$ws.Range("I74").Value = 1463.9259
$ws.Range("J74").Value = 2550.5715
$ws.Range("K74").Value = 1463.9259
$ws.Range("L74").Value = 2550.5715
$ws.Range("M74").Value = -589.9259
$ws.Range("N74").Value = -4298.5715
$ws.Range("H77").Value = 1687.6471
$ws.Range("I77").Value = 1463.9259
$ws.Range("J77").Value = 2550.5715
$ws.Range("K77").Value = 7319.6295
$ws.Range("L77").Value = 12752.8575
$ws.Range("M77").Value = -2951.6295
$ws.Range("N77").Value = -21488.8575
$ws.Range("H88").Value = 9672.416999999999
$ws.Range("J88").Value = 21405.4
$ws.Range("L88").Value = 21405.4
$ws.Range("N88").Value = -22217.4
$ws.Range("H91").Value = 9672.416999999999
$ws.Range("J91").Value = 21405.4
$ws.Range("L91").Value = 21405.4
$ws.Range("N91").Value = -24213.4
$ws.Range("H116").Value = 848.7917
$ws.Range("I116").Value = 643.65
$ws.Range("K116").Value = 643.65
$ws.Range("M116").Value = 1650.35
$ws.Range("H122").Value = 1804.6364
$ws.Range("I122").Value = 905.2857
$ws.Range("K122").Value = 2715.8571
$ws.Range("M122").Value = -265.8571000000002
$ws.Range("H132").Value = 1809.5769
$ws.Range("I132").Value = 1910.8695
$ws.Range("J132").Value = 1033
$ws.Range("K132").Value = 5732.6085
$ws.Range("L132").Value = 3099
$ws.Range("M132").Value = -3202.6085
$ws.Range("N132").Value = -8159
$ws.Range("H136").Value = 3649.6155
$ws.Range("I136").Value = 3299.6086
$ws.Range("J136").Value = 6333
$ws.Range("K136").Value = 9898.825800000001
$ws.Range("L136").Value = 18999
$ws.Range("M136").Value = -7348.825800000001
$ws.Range("N136").Value = -24099

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 848.7917
$ws.Range("I3").Value = 643.65
$ws.Range("K3").Value = 643.65
$ws.Range("M3").Value = -529.65
$ws.Range("H4").Value = 37.5
$ws.Range("I4").Value = 37.5
$ws.Range("K4").Value = 37.5
$ws.Range("M4").Value = 77.5
$ws.Range("H13").Value = 52999.5
$ws.Range("J13").Value = 52999.5
$ws.Range("L13").Value = 52999.5
$ws.Range("N13").Value = -53335.5
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 710.2
$ws.Range("J22").Value = 149.5
$ws.Range("K22").Value = 710.2
$ws.Range("L22").Value = 149.5
$ws.Range("M22").Value = -537.2
$ws.Range("N22").Value = -495.5
$ws.Range("H86").Value = 20726
$ws.Range("I86").Value = 4871.2
$ws.Range("J86").Value = 100000
$ws.Range("K86").Value = 4871.2
$ws.Range("L86").Value = 100000
$ws.Range("M86").Value = -3748.2
$ws.Range("N86").Value = -102246
$ws.Range("H89").Value = 20726
$ws.Range("I89").Value = 4871.2
$ws.Range("J89").Value = 100000
$ws.Range("K89").Value = 24356
$ws.Range("L89").Value = 500000
$ws.Range("M89").Value = -18740
$ws.Range("N89").Value = -511232
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H99").Value = 2019.25
$ws.Range("I99").Value = 1029.5
$ws.Range("K99").Value = 1029.5
$ws.Range("M99").Value = 468.5
$ws.Range("H105").Value = 3373.6785
$ws.Range("I105").Value = 3087.5454
$ws.Range("K105").Value = 3087.5454
$ws.Range("M105").Value = -1340.5454
$ws.Range("H107").Value = 1053.0435
$ws.Range("I107").Value = 945.8
$ws.Range("K107").Value = 945.8
$ws.Range("M107").Value = 974.2
$ws.Range("H134").Value = 6816.548
$ws.Range("I134").Value = 6939.8286
$ws.Range("J134").Value = 6200.143
$ws.Range("K134").Value = 20819.4858
$ws.Range("L134").Value = 18600.429
$ws.Range("M134").Value = -18284.4858
$ws.Range("N134").Value = -23670.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 154.36363
$ws.Range("I7").Value = 199.66667
$ws.Range("K7").Value = 199.66667
$ws.Range("M7").Value = -86.66667000000001
$ws.Range("H20").Value = 21000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 21000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 21000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -21472
$ws.Range("H22").Value = 1038.4286
$ws.Range("I22").Value = 860.3333
$ws.Range("J22").Value = 1172
$ws.Range("K22").Value = 860.3333
$ws.Range("L22").Value = 1172
$ws.Range("M22").Value = -510.3333
$ws.Range("N22").Value = -1872
$ws.Range("H30").Value = 21000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 21000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 21000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -21182
$ws.Range("H107").Value = 725.9524
$ws.Range("I107").Value = 637.1429000000001
$ws.Range("J107").Value = 903.5714
$ws.Range("K107").Value = 637.1429000000001
$ws.Range("L107").Value = 903.5714
$ws.Range("M107").Value = 1282.8571
$ws.Range("N107").Value = -4743.5714
$ws.Range("H122").Value = 1485.6666
$ws.Range("I122").Value = 877.5
$ws.Range("K122").Value = 2632.5
$ws.Range("M122").Value = -182.5
$ws.Range("H128").Value = 21000
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 21000
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 21000
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -30960
$ws.Range("H132").Value = 2276.6
$ws.Range("I132").Value = 2298
$ws.Range("J132").Value = 2191
$ws.Range("K132").Value = 6894
$ws.Range("L132").Value = 6573
$ws.Range("M132").Value = -4364
$ws.Range("N132").Value = -11633
$ws.Range("H134").Value = 987.5
$ws.Range("I134").Value = 1059.75
$ws.Range("J134").Value = 698.5
$ws.Range("K134").Value = 3179.25
$ws.Range("L134").Value = 2095.5
$ws.Range("M134").Value = -644.25
$ws.Range("N134").Value = -7165.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7552
$ws.Range("I3").Value = 4030
$ws.Range("K3").Value = 12090
$ws.Range("M3").Value = -11978
$ws.Range("H5").Value = 749.5
$ws.Range("I5").Value = 509.6
$ws.Range("J5").Value = 1149.3334
$ws.Range("K5").Value = 1528.8
$ws.Range("L5").Value = 3448.0002
$ws.Range("M5").Value = -1416.8
$ws.Range("N5").Value = -3672.0002
$ws.Range("H11").Value = 8474782
$ws.Range("I11").Value = 8474782
$ws.Range("K11").Value = 25424346
$ws.Range("M11").Value = -25424206
$ws.Range("H12").Value = 117.85714
$ws.Range("I12").Value = 116.666664
$ws.Range("J12").Value = 118.75
$ws.Range("K12").Value = 349.999992
$ws.Range("L12").Value = 356.25
$ws.Range("M12").Value = -176.999992
$ws.Range("N12").Value = -702.25
$ws.Range("H19").Value = 2800
$ws.Range("J19").Value = 2500
$ws.Range("L19").Value = 7500
$ws.Range("N19").Value = -7848
$ws.Range("H26").Value = 480.18182
$ws.Range("J26").Value = 641.2857
$ws.Range("L26").Value = 1923.8571
$ws.Range("N26").Value = -2499.8571
$ws.Range("H70").Value = 1702
$ws.Range("I70").Value = 1212
$ws.Range("K70").Value = 3636
$ws.Range("M70").Value = -3321
$ws.Range("H73").Value = 1702
$ws.Range("I73").Value = 1212
$ws.Range("K73").Value = 3636
$ws.Range("M73").Value = -2544
$ws.Range("H125").Value = 29500
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 29500
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 88500
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -98340
$ws.Range("H130").Value = 24977.666
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H131").Value = 15874360
$ws.Range("I131").Value = 8548103
$ws.Range("J131").Value = 22223782
$ws.Range("K131").Value = 25644309
$ws.Range("L131").Value = 66671346
$ws.Range("M131").Value = -25639269
$ws.Range("N131").Value = -66681426
$ws.Range("H135").Value = 749.5
$ws.Range("I135").Value = 509.6
$ws.Range("J135").Value = 1149.3334
$ws.Range("K135").Value = 4586.400000000001
$ws.Range("L135").Value = 10344.0006
$ws.Range("M135").Value = -2051.400000000001
$ws.Range("N135").Value = -15414.0006
$ws.Range("H137").Value = 57583.9
$ws.Range("I137").Value = 105813.4
$ws.Range("J137").Value = 9354.4
$ws.Range("K137").Value = 317440.2
$ws.Range("L137").Value = 28063.2
$ws.Range("M137").Value = -312340.2
$ws.Range("N137").Value = -38263.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H43").Value = 2900
$ws.Range("I43").Value = 2900
$ws.Range("K43").Value = 2900
$ws.Range("M43").Value = -2749
$ws.Range("H53").Value = 46999.5
$ws.Range("I53").Value = 44000
$ws.Range("J53").Value = 49999
$ws.Range("K53").Value = 44000
$ws.Range("L53").Value = 49999
$ws.Range("M53").Value = -43369
$ws.Range("N53").Value = -51261
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H102").Value = 12431.037
$ws.Range("I102").Value = 16066.167
$ws.Range("J102").Value = 5160.778
$ws.Range("K102").Value = 16066.167
$ws.Range("L102").Value = 5160.778
$ws.Range("M102").Value = -14444.167
$ws.Range("N102").Value = -8404.778
$ws.Range("H113").Value = 5676.6
$ws.Range("I113").Value = 4797
$ws.Range("J113").Value = 6263
$ws.Range("K113").Value = 4797
$ws.Range("L113").Value = 6263
$ws.Range("M113").Value = -2627
$ws.Range("N113").Value = -10603
$ws.Range("H122").Value = 3446.5
$ws.Range("I122").Value = 2975.12
$ws.Range("J122").Value = 4517.8184
$ws.Range("K122").Value = 8925.360000000001
$ws.Range("L122").Value = 13553.4552
$ws.Range("M122").Value = -6475.360000000001
$ws.Range("N122").Value = -18453.4552
$ws.Range("H126").Value = 5269.857
$ws.Range("I126").Value = 4844.0713
$ws.Range("K126").Value = 14532.2139
$ws.Range("M126").Value = -12062.2139
$ws.Range("H132").Value = 2410.739
$ws.Range("I132").Value = 2108.3333
$ws.Range("J132").Value = 3499.4
$ws.Range("K132").Value = 6324.999899999999
$ws.Range("L132").Value = 10498.2
$ws.Range("M132").Value = -3794.999899999999
$ws.Range("N132").Value = -15558.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2907.4546
$ws.Range("I7").Value = 3298.6
$ws.Range("J7").Value = 2581.5
$ws.Range("K7").Value = 3298.6
$ws.Range("L7").Value = 2581.5
$ws.Range("M7").Value = -3186.6
$ws.Range("N7").Value = -2805.5
$ws.Range("H22").Value = 1384.7646
$ws.Range("I22").Value = 1249
$ws.Range("J22").Value = 1426.5385
$ws.Range("K22").Value = 1249
$ws.Range("L22").Value = 1426.5385
$ws.Range("M22").Value = -954
$ws.Range("N22").Value = -2016.5385
$ws.Range("H27").Value = 1384.7646
$ws.Range("I27").Value = 1249
$ws.Range("J27").Value = 1426.5385
$ws.Range("K27").Value = 1249
$ws.Range("L27").Value = 1426.5385
$ws.Range("M27").Value = -1142
$ws.Range("N27").Value = -1640.5385
$ws.Range("H40").Value = 2467
$ws.Range("I40").Value = 2249.4167
$ws.Range("K40").Value = 2249.4167
$ws.Range("M40").Value = -2113.4167
$ws.Range("H43").Value = 20000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19807
$ws.Range("H46").Value = 2156.375
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2156.375
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2156.375
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2532.375
$ws.Range("H55").Value = 497.1875
$ws.Range("J55").Value = 557.5
$ws.Range("L55").Value = 557.5
$ws.Range("N55").Value = -903.5
$ws.Range("H68").Value = 2442.6667
$ws.Range("I68").Value = 2403.6365
$ws.Range("J68").Value = 2550
$ws.Range("K68").Value = 2403.6365
$ws.Range("L68").Value = 2550
$ws.Range("M68").Value = -1654.6365
$ws.Range("N68").Value = -4048
$ws.Range("H71").Value = 2442.6667
$ws.Range("I71").Value = 2403.6365
$ws.Range("J71").Value = 2550
$ws.Range("K71").Value = 12018.1825
$ws.Range("L71").Value = 12750
$ws.Range("M71").Value = -8274.182500000001
$ws.Range("N71").Value = -20238
$ws.Range("H82").Value = 2967.4
$ws.Range("I82").Value = 2424.5
$ws.Range("K82").Value = 2424.5
$ws.Range("M82").Value = -2063.5
$ws.Range("H85").Value = 2967.4
$ws.Range("I85").Value = 2424.5
$ws.Range("K85").Value = 2424.5
$ws.Range("M85").Value = -1176.5
$ws.Range("H106").Value = 19833
$ws.Range("J106").Value = 19833
$ws.Range("L106").Value = 19833
$ws.Range("N106").Value = -22357
$ws.Range("H126").Value = 2907.4546
$ws.Range("I126").Value = 3298.6
$ws.Range("J126").Value = 2581.5
$ws.Range("K126").Value = 9895.799999999999
$ws.Range("L126").Value = 7744.5
$ws.Range("M126").Value = -7425.799999999999
$ws.Range("N126").Value = -12684.5
$ws.Range("H128").Value = 89999
$ws.Range("J128").Value = 89999
$ws.Range("L128").Value = 89999
$ws.Range("N128").Value = -99959
$ws.Range("H132").Value = 2449
$ws.Range("I132").Value = 1993.5834
$ws.Range("J132").Value = 3161.8262
$ws.Range("K132").Value = 5980.7502
$ws.Range("L132").Value = 9485.4786
$ws.Range("M132").Value = -3450.7502
$ws.Range("N132").Value = -14545.4786
$ws.Range("H136").Value = 2152.1316
$ws.Range("I136").Value = 1995.1154
$ws.Range("J136").Value = 2492.3333
$ws.Range("K136").Value = 5985.3462
$ws.Range("L136").Value = 7476.999899999999
$ws.Range("M136").Value = -3435.3462
$ws.Range("N136").Value = -12576.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 2334482.8
$ws.Range("I3").Value = 3500625
$ws.Range("J3").Value = 2198
$ws.Range("K3").Value = 3500625
$ws.Range("L3").Value = 2198
$ws.Range("M3").Value = -3500511
$ws.Range("N3").Value = -2426
$ws.Range("H4").Value = 2250
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2250
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2250
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -2476
$ws.Range("H8").Value = 15999.667
$ws.Range("I8").Value = 8000
$ws.Range("J8").Value = 19999.5
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 19999.5
$ws.Range("M8").Value = -7860
$ws.Range("N8").Value = -20279.5
$ws.Range("H60").Value = 39000
$ws.Range("I60").Value = 39000
$ws.Range("K60").Value = 39000
$ws.Range("M60").Value = -38178
$ws.Range("H62").Value = 8991.5
$ws.Range("I62").Value = 8989
$ws.Range("K62").Value = 8989
$ws.Range("M62").Value = -8365
$ws.Range("H65").Value = 8991.5
$ws.Range("I65").Value = 8989
$ws.Range("K65").Value = 44945
$ws.Range("M65").Value = -41825
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 7522.625
$ws.Range("I81").Value = 13749.75
$ws.Range("J81").Value = 1295.5
$ws.Range("K81").Value = 27499.5
$ws.Range("L81").Value = 2591
$ws.Range("M81").Value = -26438.5
$ws.Range("N81").Value = -4713
$ws.Range("H84").Value = 7522.625
$ws.Range("I84").Value = 13749.75
$ws.Range("J84").Value = 1295.5
$ws.Range("K84").Value = 137497.5
$ws.Range("L84").Value = 12955
$ws.Range("M84").Value = -132193.5
$ws.Range("N84").Value = -23563
$ws.Range("H113").Value = 1173.0769
$ws.Range("I113").Value = 1106.8889
$ws.Range("K113").Value = 3320.6667
$ws.Range("M113").Value = -1150.6667
$ws.Range("H126").Value = 2201.5
$ws.Range("J126").Value = 2319.6
$ws.Range("L126").Value = 6958.799999999999
$ws.Range("N126").Value = -11898.8
$ws.Range("H132").Value = 2335.4167
$ws.Range("I132").Value = 2047.55
$ws.Range("J132").Value = 3774.75
$ws.Range("K132").Value = 6142.65
$ws.Range("L132").Value = 11324.25
$ws.Range("M132").Value = -3612.65
$ws.Range("N132").Value = -16384.25
$ws.Range("H136").Value = 2651.3635
$ws.Range("I136").Value = 2498.2273
$ws.Range("J136").Value = 2957.6365
$ws.Range("K136").Value = 7494.6819
$ws.Range("L136").Value = 8872.9095
$ws.Range("M136").Value = -4944.6819
$ws.Range("N136").Value = -13972.9095
